# Update "Main Dashboard" figures on the sole worksheet.
# Values below come from the target OOXML diff (rows are keyed by the
# sheet's own row number, matching the cell references in the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value overwrites -------------------------------------------------
$ws.Range("P2").Value  = 20000
$ws.Range("S2").Value  = 78622
$ws.Range("T2").Value  = 5444.512500000001
$ws.Range("U2").Value  = 65000

$ws.Range("S3").Value  = 75199
$ws.Range("T3").Value  = 5387.3435

$ws.Range("B4").Value  = 12500
$ws.Range("I4").Value  = 10000
$ws.Range("S4").Value  = 69820
$ws.Range("T4").Value  = 5357.2225
$ws.Range("U4").Value  = 22500

$ws.Range("B5").Value  = 12500
$ws.Range("I5").Value  = 10000
$ws.Range("S5").Value  = 69709
$ws.Range("T5").Value  = 5289.154500000001
$ws.Range("U5").Value  = 22500

$ws.Range("B6").Value  = 12500
$ws.Range("I6").Value  = 10000
$ws.Range("S6").Value  = 69148
$ws.Range("T6").Value  = 5361.8285
$ws.Range("U6").Value  = 22500

$ws.Range("B7").Value  = 12500
$ws.Range("I7").Value  = 10000
$ws.Range("S7").Value  = 71026
$ws.Range("T7").Value  = 5442.6015
$ws.Range("U7").Value  = 22500

$ws.Range("S8").Value  = 69369
$ws.Range("T8").Value  = 5816.530999999999

$ws.Range("I9").Value  = 10000
$ws.Range("S9").Value  = 83438
$ws.Range("T9").Value  = 7034.6115
$ws.Range("U9").Value  = 22500

$ws.Range("P10").Value = 20000
$ws.Range("S10").Value = 67948
$ws.Range("T10").Value = 8596.6895
$ws.Range("U10").Value = 65000

$ws.Range("P11").Value = 20000
$ws.Range("T11").Value = 13942.873
$ws.Range("U11").Value = 65000

$ws.Range("T12").Value = 15769.3095

$ws.Range("T13").Value = 15209.229
$ws.Range("T14").Value = 15159.795
$ws.Range("T15").Value = 15449.4445
$ws.Range("T16").Value = 15588.433
$ws.Range("T17").Value = 16020.6585
$ws.Range("T18").Value = 16262.3825
$ws.Range("T19").Value = 15862.609
$ws.Range("T20").Value = 14887.418
$ws.Range("T21").Value = 13341.0375
$ws.Range("T22").Value = 11837.672
$ws.Range("T23").Value = 9490.991999999998
$ws.Range("T24").Value = 6891.794

$ws.Range("P25").Value = 20000
$ws.Range("T25").Value = 5770.9715
$ws.Range("U25").Value = 65000

# --- Cells removed entirely in the target (TOTAL_SS_LOAD/WESM_RATE/CURRENT_RATE
#     no longer populated for rows 11 & 12) -----------------------------------
$ws.Range("S11").ClearContents()
$ws.Range("S12").ClearContents()
$ws.Range("V12").ClearContents()
$ws.Range("W12").ClearContents()

Write-Host "Main Dashboard updated"
